# Refactor: rebuild the data rows of the report with the new dataset.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New dataset (Item, Descrição, Concentração_Obtida, Laboratório, Registro, PDF)
$data = @(
    @(1,  "sodica;heparina",               "5000ui",                       "CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA", "102980371", "OK"),
    @(4,  "PROMETAZINA",                   "25mg",                         "LABORATÓRIO TEUTO BRASILEIRO S/A",                "103700691", "OK"),
    @(5,  "HALOPERIDOL",                   "5mg",                          "CELLERA FARMACÊUTICA S.A.",                       "112360011", "OK"),
    @(6,  "clorpromazina",                 "40mg/ml",                      "SANOFI MEDLEY FARMACÊUTICA LTDA",                 "183260385", "OK"),
    @(7,  "HALOPERIDOL",                   "2mg/ml",                       "CELLERA FARMACÊUTICA S.A.",                       "112360011", "OK"),
    @(8,  "clorpromazina",                 "25mg",                         "SANOFI MEDLEY FARMACÊUTICA LTDA",                 "183260385", "OK"),
    @(9,  "codeina",                       "30mg",                         "CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA", "102980199", "OK"),
    @(10, "imipramina",                    "25mg",                         "CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA", "102980023", "OK"),
    @(11, "RISPERIDONA",                   "3mg",                          "VITAMEDIC INDUSTRIA FARMACEUTICA LTDA",           "103920197", "OK"),
    @(13, "RISPERIDONA",                   "1mg",                          "VITAMEDIC INDUSTRIA FARMACEUTICA LTDA",           "103920197", "OK"),
    @(14, "levomepromazina",               "4%",                           "SANOFI MEDLEY FARMACÊUTICA LTDA",                 "183260316", "OK"),
    @(16, "LIDOCAÍNA",                     "2%",                           "CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA", "102980249", "OK"),
    @(17, "NITRATO DE CERIO;SULFADIAZINA", "Concentração não encontrada",  "CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA", "102980560", "OK"),
    @(18, "COLAGENASE",                    "30g",                          "CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA", "102980431", "OK")
)

# Clear out the previous data body (rows below the header) before writing the new table.
$ws.Rows("2:5").ClearContents()

$lastRow = 1 + $data.Count
$bodyRange = $ws.Range("B2:F$lastRow")

# Columns B, C, D, E, F hold free-form text (including values that look like
# numbers or percentages, e.g. registry codes and "4%"/"2%" concentrations),
# so force them to Text format to avoid Excel auto-converting them while
# writing, then restore the default "Normal" style afterwards so the cells
# don't keep a lingering explicit number format.
$bodyRange.NumberFormat = "@"

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $ws.Cells.Item($row, 6).Value = $entry[5]
    $row++
}

$bodyRange.Style = "Normal"
